$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044152094170697
$ws.Range("D2").Value = 1.049924100118101
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.059091091077023
$ws.Range("I2").Value = 1.043846791376833
$ws.Range("J2").Value = 1.049218712252233
$ws.Range("K2").Value = 1.052679806216123
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.061821541926077
$ws.Range("N2").Value = 1.020271793948582

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045056499513778
$ws.Range("D3").Value = 1.050620725571588
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.059925727354995
$ws.Range("I3").Value = 1.044072927145614
$ws.Range("J3").Value = 1.049770479454477
$ws.Range("K3").Value = 1.053189073776213
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.06247028769972
$ws.Range("N3").Value = 1.020456931824559

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045642115502671
$ws.Range("D4").Value = 1.051071872512366
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.060466513178598
$ws.Range("I4").Value = 1.044218321331843
$ws.Range("J4").Value = 1.050127264457776
$ws.Range("K4").Value = 1.053518310150175
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.06289013947644
$ws.Range("N4").Value = 1.020576595665698

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045888404255405
$ws.Range("D5").Value = 1.051261625040017
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.06069403037028
$ws.Range("I5").Value = 1.044279221608695
$ws.Range("J5").Value = 1.050277197079108
$ws.Range("K5").Value = 1.053656649693768
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.063066660567513
$ws.Range("N5").Value = 1.020626870150089

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045929762802323
$ws.Range("D6").Value = 1.051293490579048
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.060732241478439
$ws.Range("I6").Value = 1.044289433912913
$ws.Range("J6").Value = 1.050302367888373
$ws.Range("K6").Value = 1.053679873291501
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.063096300100218
$ws.Range("N6").Value = 1.020635309557269

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045645406050174
$ws.Range("D7").Value = 1.051074407641785
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.060469552605229
$ws.Range("I7").Value = 1.044219135963042
$ws.Range("J7").Value = 1.050129268101094
$ws.Range("K7").Value = 1.05352015893196
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.062892498101031
$ws.Range("N7").Value = 1.020577267562841

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044457657540466
$ws.Range("D8").Value = 1.050159447603053
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.059373010189719
$ws.Range("I8").Value = 1.04392340737562
$ws.Range("J8").Value = 1.049405234700012
$ws.Range("K8").Value = 1.052851975746088
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.062040773331432
$ws.Range("N8").Value = 1.020334389396667

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042367850972107
$ws.Range("D9").Value = 1.048550176064361
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.057446350858306
$ws.Range("I9").Value = 1.043395194124198
$ws.Range("J9").Value = 1.048127564349108
$ws.Range("K9").Value = 1.051672349825975
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.060540519220062
$ws.Range("N9").Value = 1.019905408692423

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040976842141321
$ws.Range("D10").Value = 1.047479437833974
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.056165764752929
$ws.Range("I10").Value = 1.043038313271726
$ws.Range("J10").Value = 1.047274616614653
$ws.Range("K10").Value = 1.050884516945527
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.059540827545199
$ws.Range("N10").Value = 1.019618774865737

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040375055641217
$ws.Range("D11").Value = 1.047016315381879
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.055612190946599
$ws.Range("I11").Value = 1.042882663489282
$ws.Range("J11").Value = 1.046905017152179
$ws.Range("K11").Value = 1.050543054372008
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.05910807974964
$ws.Range("N11").Value = 1.019494510989721

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040151606054811
$ws.Range("D12").Value = 1.046844369840537
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.055406710083292
$ws.Range("I12").Value = 1.042824680841677
$ws.Range("J12").Value = 1.046767692312591
$ws.Range("K12").Value = 1.050416172097807
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.058947358079062
$ws.Range("N12").Value = 1.019448331791124

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040199533078416
$ws.Range("D13").Value = 1.046881249154888
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.055450780015237
$ws.Range("I13").Value = 1.042837125881722
$ws.Range("J13").Value = 1.046797150711416
$ws.Range("K13").Value = 1.050443390925443
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.058981832491167
$ws.Range("N13").Value = 1.019458238377856

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040356583575412
$ws.Range("D14").Value = 1.047002100694514
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.055595202946704
$ws.Range("I14").Value = 1.042877874036838
$ws.Range("J14").Value = 1.046893666630577
$ws.Range("K14").Value = 1.050532567214013
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.059094794028046
$ws.Range("N14").Value = 1.019490694250484

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040453358300207
$ws.Range("D15").Value = 1.047076571797567
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.055684205426132
$ws.Range("I15").Value = 1.042902958156949
$ws.Range("J15").Value = 1.046953128118153
$ws.Range("K15").Value = 1.050587505357196
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.059164396065405
$ws.Range("N15").Value = 1.019510688473422

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041016791984542
$ws.Range("D16").Value = 1.047510184708528
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.05620252332168
$ws.Range("I16").Value = 1.043048619729395
$ws.Range("J16").Value = 1.047299140167363
$ws.Range("K16").Value = 1.050907171920151
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.059569550342355
$ws.Range("N16").Value = 1.01962701873385

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041370361631299
$ws.Range("D17").Value = 1.047782317371452
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.056527900147212
$ws.Range("I17").Value = 1.043139690470087
$ws.Range("J17").Value = 1.047516113510719
$ws.Range("K17").Value = 1.051107603853916
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.059823727382396
$ws.Range("N17").Value = 1.0196999498827

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041576643835508
$ws.Range("D18").Value = 1.047941097268782
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.056717776416447
$ws.Range("I18").Value = 1.04319270258068
$ws.Range("J18").Value = 1.047642644514663
$ws.Range("K18").Value = 1.051224480909281
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.059971996551328
$ws.Range("N18").Value = 1.01974247493263

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041646989353407
$ws.Range("D19").Value = 1.047995245464338
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.056782534444392
$ws.Range("I19").Value = 1.043210760019085
$ws.Range("J19").Value = 1.047685783896432
$ws.Range("K19").Value = 1.051264327589275
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.060022554502353
$ws.Range("N19").Value = 1.019756972411543

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04133242166961
$ws.Range("D20").Value = 1.047753114979357
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.056492981027039
$ws.Range("I20").Value = 1.043129930599104
$ws.Range("J20").Value = 1.047492836993376
$ws.Range("K20").Value = 1.051086102654825
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.059796455340676
$ws.Range("N20").Value = 1.019692126552899

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040310333880961
$ws.Range("D21").Value = 1.046966510742134
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.055552670085154
$ws.Range("I21").Value = 1.042865879336828
$ws.Range("J21").Value = 1.046865246170148
$ws.Range("K21").Value = 1.050506308332579
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.059061529104796
$ws.Range("N21").Value = 1.019481137410716

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039668174288234
$ws.Range("D22").Value = 1.046472397846075
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.054962276411633
$ws.Range("I22").Value = 1.042698891661552
$ws.Range("J22").Value = 1.04647042895853
$ws.Range("K22").Value = 1.050141492193696
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.058599568916213
$ws.Range("N22").Value = 1.019348352717805

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040008550529201
$ws.Range("D23").Value = 1.046734292680082
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.055275177275757
$ws.Range("I23").Value = 1.042787506564
$ws.Range("J23").Value = 1.046679750102374
$ws.Range("K23").Value = 1.050334913977785
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.05884445121686
$ws.Range("N23").Value = 1.019418756341294

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041349564947676
$ws.Range("D24").Value = 1.047766310130509
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.0565087591981
$ws.Range("I24").Value = 1.043134340998179
$ws.Range("J24").Value = 1.047503354729697
$ws.Range("K24").Value = 1.051095818218631
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.059808778365477
$ws.Range("N24").Value = 1.019695661623593

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042907734385645
$ws.Range("D25").Value = 1.048965846043363
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.057943766312583
$ws.Range("I25").Value = 1.04353258778639
$ws.Range("J25").Value = 1.048458083197417
$ws.Range("K25").Value = 1.051977565653718
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.060928292415779
$ws.Range("N25").Value = 1.020016426266702
